$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (Nithya Priya) updates
$ws.Range("F4").Value = 87
$ws.Range("N4").Value = 341
$ws.Range("Q4").Value = 2.99

# Row 6 (Mohammed Ibrahim Sultan) updates
$ws.Range("J6").Value = 392
$ws.Range("N6").Value = 392
$ws.Range("Q6").Value = 4.78
